# Generate Report for Handback
# The ddae4916-7eb4-4bbf-861d-476a0f304765 file has now been handed back
# and is in sync with en-US for both zh-cn and de-de locales, so the
# "Ready for handoff" status rows are updated to reflect the handback,
# along with the latest handback datetimes, and the stale error message
# is cleared out.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 6 (ddae4916-7eb4-4bbf-861d-476a0f304765.md) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E6").Value = "Handed back: in sync with en-US"
$overview.Range("F6").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 6 (ddae4916-7eb4-4bbf-861d-476a0f304765.md) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C6").Value = "Handed back: in sync with en-US"
$zhcn.Range("K6").Value = "2016-10-24 09:33:54"
$zhcn.Range("P6").Value = ""

# --- de-de sheet: row 6 (ddae4916-7eb4-4bbf-861d-476a0f304765.md) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C6").Value = "Handed back: in sync with en-US"
$dede.Range("K6").Value = "2016-10-24 09:34:11"
$dede.Range("P6").Value = ""
